$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5489456653594971
$ws.Range("B1").Value = 2.372661590576172
$ws.Range("C1").Value = 6.196498870849609
$ws.Range("D1").Value = 1.746636390686035
$ws.Range("E1").Value = 1.034395217895508
